$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.655.90"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.475.21"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.25"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.85"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0871"
$ws.Range("E10").Value = "  +11.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.02"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.111"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "2.858.18"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.69"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "2.494.59"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "41.631.12"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "0.0₃0958"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.49"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.46"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.51"
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.68"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.75"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.88"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +3.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.58"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.87"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.53"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0770"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.51"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.92"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.61"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "1.984.46"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0286"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.25"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "2.716.72"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.52"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.98"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.06"
$ws.Range("E51").Value = "  +2.36%  "
